# "Generate Report for Handoff"
# The a47bf404-... file just got handed off for localization, so its
# status moves from "In Translation" to "Ready for handoff", its
# priority moves from "ht" to "mt", and the relevant handoff timestamps
# are refreshed. The Status columns get visibly wider to fit the new,
# longer "Ready for handoff" text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
# Row 2 = 62499ccc-... (unchanged), Row 3 = a47bf404-... (updated)
$wsOverview.Range("E3").Value2 = "Ready for handoff"
$wsOverview.Range("F3").Value2 = "Ready for handoff"
$wsOverview.Range("G3").Value2 = "2016-08-31 06:16:21"

$wsOverview.Columns.Item(5).ColumnWidth = 17
$wsOverview.Columns.Item(6).ColumnWidth = 17

# --- zh-cn sheet ---
# Row 2 = 62499ccc-... (unchanged), Row 3 = a47bf404-... (updated)
$wsZhCn.Range("C3").Value2 = "Ready for handoff"
$wsZhCn.Range("E3").Value2 = "mt"
$wsZhCn.Range("H3").Value2 = "2016-08-31 06:16:17"

$wsZhCn.Columns.Item(3).ColumnWidth = 17

# --- de-de sheet ---
# Row 2 = 62499ccc-... (unchanged), Row 3 = a47bf404-... (updated)
$wsDeDe.Range("C3").Value2 = "Ready for handoff"
$wsDeDe.Range("E3").Value2 = "mt"
$wsDeDe.Range("H3").Value2 = "2016-08-31 06:16:21"

$wsDeDe.Columns.Item(3).ColumnWidth = 17
